# Work Breakdown Structure (Peer Review) - prep for hand-in
#
# The team's placeholder "*" row-of-one was expanded into five named rows
# (one per teammate) with their individual task contributions filled in,
# and the old two-column "main.class" / ".class" breakdown was collapsed
# into a single "Code" column. The now-unused columns (old J through old
# L, i.e. the 5 columns between "Code" and "Issue Tracking") were removed
# outright so "Issue Tracking" / "Overall individual contribution" slide
# left to butt up against the "Code" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 5 now-empty columns between the "Code" header (col G) and the
# "Issue Tracking" header (old col P). Deleting these whole columns shifts
# everything after them left, which is what turns old P/Q into new K/L.
$ws.Range("H1:L1").EntireColumn.Delete()

# The two sub-headers under "Code" ("main.class" / ".class") collapse into
# a single "Code" header now that there's only one column for it.
$ws.Range("G7").Value = "Code"

# Fill in the five team members and their task assignments / contribution.
$ws.Range("A8").Value = "Choo Kye Yong"
$ws.Range("G8").Value = 100

$ws.Range("A9").Value = "KEVIN LIM ERN KEE"
$ws.Range("F9").Value = 100

$ws.Range("A10").Value = "MISHRA ADITI RAKESH"
$ws.Range("C10").Value = 100

$ws.Range("A11").Value = "RAASHI SINGH"
$ws.Range("B11").Value = 100

$ws.Range("A12").Value = "TOM TANG GUAN LIANG"
$ws.Range("D12").Value = 50
$ws.Range("E12").Value = 50

# Leave the selection where the author finished editing.
$ws.Range("F10").Select()
